# Generate Report for Handback
# Populates the "bbad4a33-..." handback row (row 6) on the zh-cn and de-de
# sheets with the freshly generated handback target file, handback datetime
# and an "out of date" error message, and widens the Error Detail column.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4a64ae1ef246b50ee11df5f2c12c4fde3cc36335/e2e/bbad4a33-0532-423b-9acd-dcfbe91a200f.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c5f0e2613e1672cc0025e34aec12154d08cc2088/e2e/bbad4a33-0532-423b-9acd-dcfbe91a200f.md."

function Update-HandbackRow($ws, $orgSuffix, $targetFile, $handbackDateTime, $linkHash) {
    # I6: "Latest Target File" - becomes a hyperlink to the handback markdown
    $ws.Hyperlinks.Add(
        $ws.Range("I6"),
        "https://github.com/OpenLocalizationTestOrg/ol-test0-$orgSuffix/blob/$linkHash/e2e/bbad4a33-0532-423b-9acd-dcfbe91a200f.md",
        "",
        "",
        "bbad4a33-0532-423b-9acd-dcfbe91a200f.md"
    ) | Out-Null
    $ws.Range("I6").Font.Underline = 2
    $ws.Range("I6").Font.Color = 15570276

    # J6: "Latest Handback File" - the generated xliff file name
    $ws.Range("J6").Value = $targetFile

    # K6: "Latest Handback DateTime"
    $ws.Range("K6").Value = $handbackDateTime

    # P6: "Error Detail"
    $ws.Range("P6").Value = $errorDetail

    # Widen the Error Detail column so the message is readable
    $ws.Columns.Item(16).ColumnWidth = 39.16666667
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-HandbackRow $wsZhCn "zhcn" `
    "bbad4a33-0532-423b-9acd-dcfbe91a200f.8d81cd1d27ead922cf65c8b4ba0692818b223b5a.zh-cn.xlf" `
    "2016-11-03 19:34:32" `
    "259b5f08a62d8c8f52e134c8f421d8c82f4372a8"

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-HandbackRow $wsDeDe "dede" `
    "bbad4a33-0532-423b-9acd-dcfbe91a200f.8d81cd1d27ead922cf65c8b4ba0692818b223b5a.de-de.xlf" `
    "2016-11-03 19:34:49" `
    "7a1b0aa065c17531d74f001225f190807133c7b2"
